$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-10-05 01:49:35"

# Hyperlinks in this engine cannot be edited/deleted individually in place
# (Hyperlink.Delete() on a single item is a no-op, and Hyperlinks.Add()
# always appends rather than replacing an existing ref). So wipe the whole
# collection up front -- while the sheet still matches the original layout
# 1:1 -- and rebuild it from scratch once every row is in its final place.
$ws.Hyperlinks.Delete()

# Insert a new row at position 6; this pushes the old rows 6 and 7 down to
# rows 7 and 8 respectively (cell values/styles shift automatically).
$ws.Rows.Item(6).Insert()

# Refresh the "taken at" timestamp for every data row, including the new one.
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Populate the newly inserted row 6 with the new job listing.
$ws.Cells.Item(6, 2).Value = "【SalesIQ活用】CRMと連携したリード獲得方法を教えてください"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "~ 5,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5400402"
$ws.Cells.Item(6, 7).Value = 10

# Rebuild every hyperlink (F2:F8) against the now-final row layout, in order,
# so the relationship ids line up the same way Excel would regenerate them.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5406694")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5406904")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5406636")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5400402")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5406717")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5406440")
